{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// This applies three proofing-related textual fixes inside the FAQ body\n// (each originally wrapped in a <w:proofErr .../> spellcheck/grammar-check\n// marker pair, which Word removes once the flagged text is accepted/fixed):\n//   1. \"its\"   -> \"it\\u2019s\"   (smart apostrophe)   [[w:proofErr type=spellStart/spellEnd]]\n//   2. \"In\"    -> \"in\"                                 [[w:proofErr type=gramStart/gramEnd]]\n//   3. \"few..\" -> \"few:\"                                [[w:proofErr type=gramStart/gramEnd]]\n//\n// The Office.js object model has no direct handle on <w:proofErr/> marks\n// (they are not part of the Range/Paragraph content model), so each whole\n// paragraph's run sequence is rebuilt explicitly via insertOoxml \u2014 this both\n// rewrites the target word and drops the now-stale proofErr wrapper, while\n// every other run/paragraph property (paraId, rsids, pPr, rPr, etc.) is\n// round-tripped unchanged.\n\nconst body = context.document.body;\n\nasync function rebuildParagraph(searchText, buildRuns) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not locate text: \" + searchText);\n  }\n  const hitParagraphs = results.items[0].paragraphs;\n  hitParagraphs.load(\"items\");\n  await context.sync();\n  const paragraph = hitParagraphs.items[0];\n\n  const paraRange = paragraph.getRange();\n  const ooxmlResult = paraRange.getOoxml();\n  await context.sync();\n  const xml = ooxmlResult.value;\n\n  // Pull the paragraph's own attributes (w14:paraId, rsids, ...) and its\n  // <w:pPr> block straight out of the round-tripped XML so they are\n  // preserved byte-for-byte; only the run content is replaced.\n  const pMatch = xml.match(/<w:p\\b([^>]*)>/);\n  const pAttrs = pMatch ? pMatch[1] : \"\";\n  const pPrMatch = xml.match(/<w:pPr>[\\s\\S]*?<\\/w:pPr>/);\n  const pPr = pPrMatch ? pPrMatch[0] : \"\";\n\n  const runsXml = buildRuns();\n\n  const newParagraph =\n    \"<w:p\" + pAttrs + \">\" + pPr + runsXml + \"</w:p>\";\n\n  const packageXml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\" xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    newParagraph +\n    \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\";\n\n  paraRange.insertOoxml(packageXml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst RPR = '<w:rPr><w:lang w:val=\"en-US\"/></w:rPr>';\n\n// 1) \"its\" -> \"it's\" (curly apostrophe), drop spellStart/spellEnd proofErr.\nawait rebuildParagraph(\"There is no specific way to time the market, but its about\", () => {\n  return (\n    \"<w:r>\" + RPR + '<w:t xml:space=\"preserve\">There is no specific way to time the market, but </w:t></w:r>' +\n    \"<w:r>\" + RPR + \"<w:t>it\\u2019s</w:t></w:r>\" +\n    \"<w:r>\" + RPR + '<w:t xml:space=\"preserve\"> about the TIME IN the market. There are ways to understand how the market works but no way of timing it. </w:t></w:r>'\n  );\n});\n\n// 2) \"In\" -> \"in\", drop gramStart/gramEnd proofErr.\nawait rebuildParagraph(\"There are chances In most situations\", () => {\n  return (\n    \"<w:r>\" + RPR + '<w:t xml:space=\"preserve\">There are chances </w:t></w:r>' +\n    \"<w:r>\" + RPR + \"<w:t>in</w:t></w:r>\" +\n    \"<w:r>\" + RPR + '<w:t xml:space=\"preserve\"> most situations this probably won\\u2019t occur. However, we do recommend our clients about any worst-case situation to see if they are comfortable. Thus, we would also have a plan if an investment doesn\\u2019t go as accordingly. </w:t></w:r>'\n  );\n});\n\n// 3) \"few..\" -> \"few:\", drop gramStart/gramEnd proofErr.\nawait rebuildParagraph(\"There are multiple ways but here are a few..\", () => {\n  return (\n    \"<w:r>\" + RPR + '<w:t xml:space=\"preserve\">There are multiple ways but here are a </w:t></w:r>' +\n    \"<w:r>\" + RPR + \"<w:t>few:</w:t></w:r>\"\n  );\n});\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Applies three proofing-related textual fixes inside the FAQ body (each\n# originally wrapped in a <w:proofErr .../> spellcheck/grammar-check marker\n# pair, which Word removes once the flagged text is accepted/fixed):\n#   1. \"its\"   -> \"it's\"   (smart/curly apostrophe)  [w:proofErr type=spellStart/spellEnd]\n#   2. \"In\"    -> \"in\"                                [w:proofErr type=gramStart/gramEnd]\n#   3. \"few..\" -> \"few:\"                              [w:proofErr type=gramStart/gramEnd]\n#\n# The Word object model has no direct handle on <w:proofErr/> marks (they\n# are not exposed as Range/Paragraph content), so each whole paragraph's\n# run sequence is rebuilt explicitly via Range.InsertXML \u2014 this both\n# rewrites the target word and drops the now-stale proofErr wrapper, while\n# every other paragraph property (paraId, rsids, pPr, rPr, etc.) is\n# round-tripped unchanged by reusing the paragraph's own WordOpenXML.\n\nfunction Find-ParagraphIndex($doc, $searchText) {\n    $rng = $doc.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($searchText, $false, $true) | Out-Null\n    if (-not $rng.Find.Found) {\n        throw \"Could not locate text: $searchText\"\n    }\n    $paras = $doc.Paragraphs\n    for ($i = 1; $i -le $paras.Count; $i++) {\n        $p = $paras.Item($i)\n        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {\n            return $i\n        }\n    }\n    throw \"Could not locate paragraph containing: $searchText\"\n}\n\nfunction Rebuild-Paragraph($doc, $searchText, $runsXml) {\n    $idx = Find-ParagraphIndex $doc $searchText\n    $p = $doc.Paragraphs.Item($idx)\n    $pRange = $p.Range\n    $xml = $pRange.WordOpenXML\n\n    # Pull the paragraph's own attributes (w14:paraId, rsids, ...) and its\n    # <w:pPr> block straight out of the round-tripped XML so they are\n    # preserved byte-for-byte; only the run content is replaced.\n    if ($xml -match '<w:p\\b([^>]*)>') {\n        $pAttrs = $matches[1]\n    } else {\n        $pAttrs = \"\"\n    }\n    if ($xml -match '<w:pPr>[\\s\\S]*?</w:pPr>') {\n        $pPr = $matches[0]\n    } else {\n        $pPr = \"\"\n    }\n\n    $newParagraph = \"<w:p\" + $pAttrs + \">\" + $pPr + $runsXml + \"</w:p>\"\n\n    $packageXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n        '<pkg:xmlData>' + `\n        '<w:document xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\" xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n        '<w:body>' + $newParagraph + '</w:body></w:document>' + `\n        '</pkg:xmlData></pkg:part></pkg:package>'\n\n    $pRange.InsertXML($packageXml)\n}\n\n$d = $word.ActiveDocument\n$q = [char]0x2019\n$RPR = '<w:rPr><w:lang w:val=\"en-US\"/></w:rPr>'\n\n# 1) \"its\" -> \"it's\" (curly apostrophe), drop spellStart/spellEnd proofErr.\n$runs1 = \"<w:r>$RPR\" + '<w:t xml:space=\"preserve\">There is no specific way to time the market, but </w:t></w:r>' + `\n         \"<w:r>$RPR<w:t>it${q}s</w:t></w:r>\" + `\n         \"<w:r>$RPR\" + '<w:t xml:space=\"preserve\"> about the TIME IN the market. There are ways to understand how the market works but no way of timing it. </w:t></w:r>'\nRebuild-Paragraph $d \"There is no specific way to time the market, but its about\" $runs1\n\n# 2) \"In\" -> \"in\", drop gramStart/gramEnd proofErr.\n$runs2 = \"<w:r>$RPR\" + '<w:t xml:space=\"preserve\">There are chances </w:t></w:r>' + `\n         \"<w:r>$RPR<w:t>in</w:t></w:r>\" + `\n         \"<w:r>$RPR\" + \"<w:t xml:space=`\"preserve`\"> most situations this probably won${q}t occur. However, we do recommend our clients about any worst-case situation to see if they are comfortable. Thus, we would also have a plan if an investment doesn${q}t go as accordingly. </w:t></w:r>\"\nRebuild-Paragraph $d \"There are chances In most situations\" $runs2\n\n# 3) \"few..\" -> \"few:\", drop gramStart/gramEnd proofErr.\n$runs3 = \"<w:r>$RPR\" + '<w:t xml:space=\"preserve\">There are multiple ways but here are a </w:t></w:r>' + `\n         \"<w:r>$RPR<w:t>few:</w:t></w:r>\"\nRebuild-Paragraph $d \"There are multiple ways but here are a few..\" $runs3\n\nWrite-Output \"done\"\n"}
